$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Commit: "homogenization of data cohort and add FED3 data"
# On Sheet1 the two slope columns ("BW_SLOPE_LM" in D and
# "BW_TIME_INDIVIDUAL_COEFF_LMER" in E) were removed; the surviving
# "cumsumFI_SLOPE_LM"/"cumsumFI_TIME_INDIVIDUAL_COEFF_LMER" columns
# (old F:G) shift left into D:E. Reproduce this with a plain column
# delete (select first, so the resulting selection/active cell lands on
# the new D:E columns exactly like it would after a manual selection +
# delete in the Excel UI).
[void]$ws.Columns("D:E").Select()
$ws.Columns("D:E").Delete()
